$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Jan-2024")

# --- Row 1: extend date series into S1/T1 ---
$ws.Range("S1:T1").NumberFormat = $ws.Range("R1").NumberFormat
$ws.Range("S1").Value = 45307
$ws.Range("T1").Value = 45308

# --- Row 2: H2:T2 = Absent ---
$ws.Range("H2:T2").Value = "Absent"

# --- Row 3: H3:K3 = Present, L3:T3 = Absent ---
$ws.Range("H3:K3").Value = "Present"
$ws.Range("L3:T3").Value = "Absent"

# --- Row 4: H4:K4 = Present, L4:T4 = Absent ---
$ws.Range("H4:K4").Value = "Present"
$ws.Range("L4:T4").Value = "Absent"

# --- Row 5: H5:K5 = Present, L5:T5 = Absent ---
$ws.Range("H5:K5").Value = "Present"
$ws.Range("L5:T5").Value = "Absent"

# --- Comments H2:N2, matching the existing G2 "university Exam" note ---
foreach ($col in @("H", "I", "J", "K", "L", "M", "N")) {
    $cmt = $ws.Range($col + "2").AddComment()
    $cmt.Text("A:`r`nuniversity Exam")
}

# --- Data validation: extend list validation from C2:G5 to C2:T5 ---
$ws.Range("C2:T5").Validation.Delete()
$ws.Range("C2:T5").Validation.Add(3, 1, 1, '"Present, Absent,Reason"')

# --- Sheet view: zoom to 70%, change selection, drop topLeftCell override ---
$excel.ActiveWindow.Zoom = 70
$ws.Range("R13").Select() | Out-Null
